$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = ""

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = ""

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = ""

# Row 6 - DÜZENLİ EFT
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F8").Value = ""

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F9").Value = ""

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F10").Value = ""

# Row 11 - DÜZENLİ HAVALE
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("F13").Value = ""

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("F14").Value = ""
